# Update template version string 1.0.2 -> 1.0.3
$wb = $excel.ActiveWorkbook
$isaSheet = $wb.Worksheets.Item("isa_template")
$isaSheet.Range("B4").Value = "1.0.3"

# Remove the "Protocol REF" column (column E) from the annotation table on
# the "New Table" sheet. Deleting the whole worksheet column shifts every
# later column one place to the left (and keeps the row data intact).
$ws = $wb.Worksheets.Item("New Table")
$ws.Columns("E:E").Delete()

# Shrink the table range to match the new (one narrower) extent.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:W2"))

# Resizing keeps the old column headers/count mismatched with the sheet
# (it just truncates the last column instead of dropping the deleted one),
# so fix up the header row -- and therefore the table's column names --
# by writing the correct header text back into each header cell from the
# deleted column onward.
$ws.Range("E1").Value = "Parameter [solvent extraction]"
$ws.Range("F1").Value = "Term Source REF (CHMO:0001598)"
$ws.Range("G1").Value = "Term Accession Number (CHMO:0001598)"
$ws.Range("H1").Value = "Parameter [derivatisation]"
$ws.Range("I1").Value = "Term Source REF (CHMO:0001485)"
$ws.Range("J1").Value = "Term Accession Number (CHMO:0001485)"
$ws.Range("K1").Value = "Parameter [chromatography instrument]"
$ws.Range("L1").Value = "Term Source REF (OBI:0000485)"
$ws.Range("M1").Value = "Term Accession Number (OBI:0000485)"
$ws.Range("N1").Value = "Parameter [chromatography column model]"
$ws.Range("O1").Value = "Term Source REF (DPBO:0000048)"
$ws.Range("P1").Value = "Term Accession Number (DPBO:0000048)"
$ws.Range("Q1").Value = "Parameter [chromatography column type]"
$ws.Range("R1").Value = "Term Source REF (DPBO:0000053)"
$ws.Range("S1").Value = "Term Accession Number (DPBO:0000053)"
$ws.Range("T1").Value = "Parameter [label]"
$ws.Range("U1").Value = "Term Source REF (CHEBI:35209)"
$ws.Range("V1").Value = "Term Accession Number (CHEBI:35209)"
$ws.Range("W1").Value = "Output [Raw Data File]"
